# Fix Training Data Issue (#48)
# The "Date" column (BF) held a mangled string ("6-26-2013-14") that was
# actually supposed to be the ISO date 2014-06-26 (NBA stats for the day
# were off by one, per the commit message). Correct every data row.
#
# The cells must stay plain text (not get silently reinterpreted as a
# date serial number by Excel's normal "smart" entry), so the column is
# pre-formatted as Text ("@") before the literal value is written in,
# exactly like a user would do in the Excel UI to force literal text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 2
$lastDataRow = 31
$dateCol = 58   # column BF

$dateRange = $ws.Range($ws.Cells.Item($firstDataRow, $dateCol), $ws.Cells.Item($lastDataRow, $dateCol))
$dateRange.NumberFormat = "@"

for ($row = $firstDataRow; $row -le $lastDataRow; $row++) {
    $ws.Cells.Item($row, $dateCol).Value = "2014-06-26"
}
